$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 12: effort values changed ---
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 2

# --- Row 18: the 2.25h entry moved from column C to column B ---
$ws.Range("B18").Value = 2.25
$ws.Range("C18").Value = $null

# --- New rows 21-23 ---
$ws.Range("A21").Value = 41452
$ws.Range("B21").Value = 2.5
$ws.Range("D21").Value = "TODOs, code cleanup"

$ws.Range("A22").Value = 41455
$ws.Range("B22").Value = 1.5
$ws.Range("D22").Value = "Concept of tc14"

$ws.Range("A23").Value = 41456
$ws.Range("B23").Value = 2
$ws.Range("D23").Value = "Implementation tc14"

# Copy the date formatting (style index used by column A, e.g. row 19) onto
# the new date cells without disturbing their values.
$ws.Range("A19").Copy() | Out-Null
$ws.Range("A21:A23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Selection moves on to the next empty row, as Excel would leave it
#     after the user finished typing the new entries ---
$ws.Range("A24").Select() | Out-Null
